$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 60.083332
$ws.Range("I38").Value = 60.083332
$ws.Range("K38").Value = 180.249996
$ws.Range("M38").Value = 191.750004
# Row 40
$ws.Range("H40").Value = 1321.6721
$ws.Range("I40").Value = 1128.1364
$ws.Range("J40").Value = 1430.8462
$ws.Range("K40").Value = 1128.1364
$ws.Range("L40").Value = 1430.8462
$ws.Range("M40").Value = -953.1364000000001
$ws.Range("N40").Value = -1780.8462
# Row 43
$ws.Range("H43").Value = 21513.8
$ws.Range("I43").Value = 1174
$ws.Range("K43").Value = 1174
$ws.Range("M43").Value = -1105
# Row 58
$ws.Range("H58").Value = 397.94446
# Row 100
$ws.Range("H100").Value = 1563.5714
$ws.Range("I100").Value = 1582.5
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 1582.5
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -1041.5
$ws.Range("N100").Value = -2532
# Row 113
$ws.Range("H113").Value = 56638.047
$ws.Range("I113").Value = 92249
$ws.Range("J113").Value = 13904.9
$ws.Range("K113").Value = 92249
$ws.Range("L113").Value = 13904.9
$ws.Range("M113").Value = -88995
$ws.Range("N113").Value = -20412.9
# Row 135
$ws.Range("H135").Value = 2311.5293
$ws.Range("I135").Value = 1731.4
$ws.Range("J135").Value = 3140.2856
$ws.Range("K135").Value = 15582.6
$ws.Range("L135").Value = 28262.5704
$ws.Range("M135").Value = -13047.6
$ws.Range("N135").Value = -33332.5704

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 8824.8125
$ws.Range("I45").Value = 10957.909
$ws.Range("K45").Value = 10957.909
$ws.Range("M45").Value = -10580.909
# Row 74
$ws.Range("H74").Value = 10158.9
$ws.Range("I74").Value = 2198.5
$ws.Range("J74").Value = 15465.833
$ws.Range("K74").Value = 2198.5
$ws.Range("L74").Value = 15465.833
$ws.Range("M74").Value = -1324.5
$ws.Range("N74").Value = -17213.833
# Row 77
$ws.Range("H77").Value = 10158.9
$ws.Range("I77").Value = 2198.5
$ws.Range("J77").Value = 15465.833
$ws.Range("K77").Value = 10992.5
$ws.Range("L77").Value = 77329.16500000001
$ws.Range("M77").Value = -6624.5
$ws.Range("N77").Value = -86065.16500000001
# Row 122
$ws.Range("H122").Value = 1976.7894
$ws.Range("I122").Value = 1603.6428
$ws.Range("K122").Value = 4810.928400000001
$ws.Range("M122").Value = -2360.928400000001
# Row 125
$ws.Range("H125").Value = 35186.625
$ws.Range("J125").Value = 35186.625
$ws.Range("L125").Value = 35186.625
$ws.Range("N125").Value = -45026.625
# Row 132
$ws.Range("H132").Value = 3080.4285
$ws.Range("I132").Value = 2619.5
$ws.Range("J132").Value = 4002.2856
$ws.Range("K132").Value = 7858.5
$ws.Range("L132").Value = 12006.8568
$ws.Range("M132").Value = -5328.5
$ws.Range("N132").Value = -17066.8568
# Row 139
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6168.5
$ws.Range("I86").Value = 4299.6665
$ws.Range("J86").Value = 8037.3335
$ws.Range("K86").Value = 4299.6665
$ws.Range("L86").Value = 8037.3335
$ws.Range("M86").Value = -3176.6665
$ws.Range("N86").Value = -10283.3335
# Row 89
$ws.Range("H89").Value = 6168.5
$ws.Range("I89").Value = 4299.6665
$ws.Range("J89").Value = 8037.3335
$ws.Range("K89").Value = 21498.3325
$ws.Range("L89").Value = 40186.6675
$ws.Range("M89").Value = -15882.3325
$ws.Range("N89").Value = -51418.6675

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 515789
$ws.Range("J9").Value = 515789
$ws.Range("L9").Value = 515789
$ws.Range("N9").Value = -516125
# Row 33
$ws.Range("H33").Value = 36371.125
$ws.Range("I33").Value = 19000
$ws.Range("J33").Value = 42161.5
$ws.Range("K33").Value = 19000
$ws.Range("L33").Value = 42161.5
$ws.Range("M33").Value = -18621
$ws.Range("N33").Value = -42919.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
# Row 139
$ws.Range("H139").Value = 3403.1667
$ws.Range("I139").Value = 3403.1667
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 10209.5001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -5069.500100000001
$ws.Range("N139").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 256.2857
$ws.Range("I2").Value = 255.79167
$ws.Range("J2").Value = 257.36365
$ws.Range("K2").Value = 255.79167
$ws.Range("L2").Value = 257.36365
$ws.Range("M2").Value = -142.79167
$ws.Range("N2").Value = -483.36365
# Row 62
$ws.Range("H62").Value = 46000
$ws.Range("I62").Value = 46000
$ws.Range("K62").Value = 46000
$ws.Range("M62").Value = -45314
# Row 65
$ws.Range("H65").Value = 46000
$ws.Range("I65").Value = 46000
$ws.Range("K65").Value = 138000
$ws.Range("M65").Value = -134568
# Row 102
$ws.Range("H102").Value = 5915.8335
$ws.Range("I102").Value = 5599
$ws.Range("K102").Value = 5599
$ws.Range("M102").Value = -3977
# Row 113
$ws.Range("H113").Value = 1997.5
$ws.Range("I113").Value = 1997.5
$ws.Range("K113").Value = 1997.5
$ws.Range("M113").Value = 172.5
# Row 123
$ws.Range("H123").Value = 35259
$ws.Range("J123").Value = 35259
$ws.Range("L123").Value = 35259
$ws.Range("N123").Value = -40159

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 11622.091
$ws.Range("I7").Value = 12204.889
$ws.Range("J7").Value = 8999.5
$ws.Range("K7").Value = 12204.889
$ws.Range("L7").Value = 8999.5
$ws.Range("M7").Value = -12092.889
$ws.Range("N7").Value = -9223.5
# Row 22
$ws.Range("H22").Value = 1138
$ws.Range("I22").Value = 574
$ws.Range("J22").Value = 1363.6
$ws.Range("K22").Value = 574
$ws.Range("L22").Value = 1363.6
$ws.Range("M22").Value = -279
$ws.Range("N22").Value = -1953.6
# Row 27
$ws.Range("H27").Value = 1138
$ws.Range("I27").Value = 574
$ws.Range("J27").Value = 1363.6
$ws.Range("K27").Value = 574
$ws.Range("L27").Value = 1363.6
$ws.Range("M27").Value = -467
$ws.Range("N27").Value = -1577.6
# Row 40
$ws.Range("H40").Value = 3278.261
$ws.Range("I40").Value = 2800.5
$ws.Range("K40").Value = 2800.5
$ws.Range("M40").Value = -2664.5
# Row 46
$ws.Range("H46").Value = 2029.7
$ws.Range("I46").Value = 1737.125
$ws.Range("K46").Value = 1737.125
$ws.Range("M46").Value = -1549.125
# Row 122
$ws.Range("H122").Value = 4460.1113
$ws.Range("I122").Value = 3885.4666
$ws.Range("K122").Value = 11656.3998
$ws.Range("M122").Value = -9206.399800000001
# Row 126
$ws.Range("H126").Value = 11622.091
$ws.Range("I126").Value = 12204.889
$ws.Range("J126").Value = 8999.5
$ws.Range("K126").Value = 36614.667
$ws.Range("L126").Value = 26998.5
$ws.Range("M126").Value = -34144.667
$ws.Range("N126").Value = -31938.5
# Row 136
$ws.Range("H136").Value = 6133.476
$ws.Range("I136").Value = 5832.385
$ws.Range("J136").Value = 6622.75
$ws.Range("K136").Value = 17497.155
$ws.Range("L136").Value = 19868.25
$ws.Range("M136").Value = -14947.155
$ws.Range("N136").Value = -24968.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 141
$ws.Range("H141").Value = 136713.33
$ws.Range("J141").Value = 136713.33
$ws.Range("L141").Value = 136713.33
$ws.Range("N141").Value = -147073.33
